# Add a new activity log row (row 22) to the "Activity Log" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

$ws.Range("A22").Value = "vaishali.kh2310@gmail.com"
$ws.Range("B22").Value = "Login"
$ws.Range("C22").Value = "2025-06-16 18:56:52"
